$d = $word.ActiveDocument

# --- Change 1: merge runs around "Tzaquitzal" (spell-checked name) into one run ---
$old1 = "Identificarse como Ludwing Juan Homero Pérez Tzaquitzal, de 25 años de edad, de nacionalidad guatemalteca, identificado con el número de DPI 2333 75953 0801, con residencia en el Cantón Chotacaj, C-92, Totonicapán, Totonicapán, Guatemala C.A. Actualmente estudiante de la Universidad Rafael Landívar Campus de Quetzaltenango en la Facultad de Ingeniería, inscrito en la carrera de Ingeniería en Informática y Sistemas, identificado con carné no. 1520909."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- Change 2: merge runs around "Chotacaj" (second occurrence) into one run ---
$old2 = "Ser una organización comunitaria, sin ánimo de lucro, sin afiliaciones políticas y religiosas con sede en el Cantón Chotacaj, C-1, Municipio de Totonicapán, Departamento de Totonicapán, Guatemala C.A."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- Change 3: "4 meses" -> "6 meses", splitting the run and inserting a _GoBack bookmark ---
$rng = $d.Content
$old3 = " cuenta con 4"
$new3 = " cuenta con 6"
$rng.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null
$bmRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Change 4: footer page number field cached result "4" -> "3" ---
$sec = $d.Sections.First
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2) | Out-Null
